$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column for rows 2-6 from 2023-11-03 (45233) to 2023-11-13 (45243)
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45243
}
